$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue 'D2' '25.777.74'
Set-TextValue 'E2' '  +0.29%  '
Set-TextValue 'D3' '1.747.26'
Set-TextValue 'E3' '  +0.07%  '
Set-TextValue 'D4' '1.002'
Set-TextValue 'E4' '  +0.01%  '
Set-TextValue 'D5' '235.53'
Set-TextValue 'E5' '  -0.32%  '
Set-TextValue 'E6' '  +0.00%  '
Set-TextValue 'E7' '  +3.33%  '
Set-TextValue 'B8' 'Cardano'
Set-TextValue 'C8' 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextValue 'D8' '0.2668'
Set-TextValue 'E8' '  +6.38%  '
Set-TextValue 'B9' 'Dogecoin'
Set-TextValue 'C9' 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextValue 'D9' '0.06187'
Set-TextValue 'E9' '  +3.05%  '
Set-TextValue 'B10' 'WrappedEther'
Set-TextValue 'C10' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D10' '1.749.83'
Set-TextValue 'E10' '  +0.24%  '
Set-TextValue 'B11' 'TRON'
Set-TextValue 'C11' 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue 'D11' '0.06935'
Set-TextValue 'E11' '  +1.43%  '
Set-TextValue 'B12' 'Solana'
Set-TextValue 'C12' 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextValue 'D12' '15.35'
Set-TextValue 'E12' '  +3.61%  '
Set-TextValue 'B13' 'Polygon'
Set-TextValue 'C13' 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue 'D13' '0.6204'
Set-TextValue 'E13' '  +10.18%  '
Set-TextValue 'B14' 'Polkadot'
Set-TextValue 'C14' 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 'D14' '4.472'
Set-TextValue 'E14' '  +0.67%  '
Set-TextValue 'B15' 'Litecoin'
Set-TextValue 'C15' 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue 'D15' '77.71'
Set-TextValue 'E15' '  +1.08%  '
Set-TextValue 'B16' 'BinanceUSD'
Set-TextValue 'C16' 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue 'D16' '1.001'
Set-TextValue 'E16' '  +0.03%  '
Set-TextValue 'B17' 'Dai'
Set-TextValue 'C17' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 'D17' '1.000'
Set-TextValue 'E17' '  -0.08%  '
Set-TextValue 'B18' 'WrappedBTC'
Set-TextValue 'C18' 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue 'D18' '25.802.40'
Set-TextValue 'E18' '  +0.24%  '
Set-TextValue 'B19' 'Avalanche'
Set-TextValue 'C19' 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue 'D19' '11.58'
Set-TextValue 'E19' '  +2.91%  '
Set-TextValue 'B20' 'ShibaInu'
Set-TextValue 'C20' 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue 'D20' '0.000006643'
Set-TextValue 'E20' '  +1.37%  '
Set-TextValue 'B21' 'WrappedliquidstakedEther2.0'
Set-TextValue 'C21' 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue 'D21' '1.975.26'
Set-TextValue 'E21' '  +0.53%  '
Set-TextValue 'B22' 'Uniswap'
Set-TextValue 'C22' 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue 'D22' '4.048'
Set-TextValue 'E22' '  +1.47%  '
Set-TextValue 'B23' 'Cosmos'
Set-TextValue 'C23' 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue 'D23' '8.246'
Set-TextValue 'E23' '  +4.92%  '
Set-TextValue 'B24' 'Chainlink'
Set-TextValue 'C24' 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 'D24' '5.136'
Set-TextValue 'E24' '  +2.65%  '
Set-TextValue 'B25' 'Monero'
Set-TextValue 'C25' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D25' '136.52'
Set-TextValue 'E25' '  -0.13%  '
Set-TextValue 'B26' 'Toncoin'
Set-TextValue 'C26' 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue 'D26' '1.460'
Set-TextValue 'E26' '  -1.54%  '
Set-TextValue 'B27' 'EthereumClassic'
Set-TextValue 'C27' 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue 'D27' '15.03'
Set-TextValue 'E27' '  +2.73%  '
Set-TextValue 'B28' 'LidoDAOToken'
Set-TextValue 'C28' 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue 'D28' '1.768'
Set-TextValue 'E28' '  -1.95%  '
Set-TextValue 'B29' 'BitcoinCash'
Set-TextValue 'C29' 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue 'D29' '102.53'
Set-TextValue 'E29' '  +0.67%  '
Set-TextValue 'B30' 'Stellar'
Set-TextValue 'C30' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 'D30' '0.08150'
Set-TextValue 'E30' '  +2.06%  '
Set-TextValue 'B31' 'InternetComputer(DFINITY)'
Set-TextValue 'C31' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 'D31' '3.692'
Set-TextValue 'E31' '  -1.67%  '
Set-TextValue 'B32' 'Filecoin'
Set-TextValue 'C32' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D32' '3.383'
Set-TextValue 'E32' '  +0.47%  '
Set-TextValue 'B33' 'Hedera'
Set-TextValue 'C33' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D33' '0.04412'
Set-TextValue 'E33' '  +0.57%  '
Set-TextValue 'B34' 'HuobiToken'
Set-TextValue 'C34' 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue 'D34' '2.648'
Set-TextValue 'E34' '  +0.42%  '
Set-TextValue 'B35' 'ARBITRUM'
Set-TextValue 'C35' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue 'D35' '0.9935'
Set-TextValue 'E35' '  +2.60%  '
Set-TextValue 'B36' 'ImmutableX'
Set-TextValue 'C36' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D36' '0.6007'
Set-TextValue 'E36' '  -0.21%  '
Set-TextValue 'B37' 'MXToken'
Set-TextValue 'C37' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 'D37' '2.627'
Set-TextValue 'E37' '  -1.95%  '
Set-TextValue 'B38' 'VeChain'
Set-TextValue 'C38' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D38' '0.01555'
Set-TextValue 'E38' '  +3.72%  '
Set-TextValue 'B39' 'RenderToken'
Set-TextValue 'C39' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D39' '1.941'
Set-TextValue 'E39' '  -3.00%  '
Set-TextValue 'B40' 'PaxDollar'
Set-TextValue 'C40' 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue 'D40' '1.001'
Set-TextValue 'E40' '  -0.06%  '
Set-TextValue 'B41' 'Quant'
Set-TextValue 'C41' 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue 'D41' '101.49'
Set-TextValue 'E41' '  -1.12%  '
Set-TextValue 'B42' 'TheSandbox'
Set-TextValue 'C42' 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue 'D42' '0.3821'
Set-TextValue 'E42' '  +2.85%  '
Set-TextValue 'B43' 'TrustWalletToken'
Set-TextValue 'C43' 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 'D43' '0.7455'
Set-TextValue 'E43' '  -1.22%  '
Set-TextValue 'B44' 'FraxShare'
Set-TextValue 'C44' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 'D44' '4.886'
Set-TextValue 'E44' '  -5.27%  '
Set-TextValue 'B45' 'Cronos'
Set-TextValue 'C45' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D45' '0.05499'
Set-TextValue 'E45' '  +4.62%  '
Set-TextValue 'B46' 'Algorand'
Set-TextValue 'C46' 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue 'D46' '0.1092'
Set-TextValue 'E46' '  +2.66%  '
Set-TextValue 'B47' 'Aptos'
Set-TextValue 'C47' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 'D47' '5.914'
Set-TextValue 'E47' '  +0.83%  '
Set-TextValue 'B48' 'Elrond'
Set-TextValue 'C48' 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-TextValue 'D48' '30.01'
Set-TextValue 'E48' '  +0.03%  '
Set-TextValue 'B49' 'Aave'
Set-TextValue 'C49' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D49' '52.58'
Set-TextValue 'E49' '  +0.70%  '
Set-TextValue 'B50' 'USDD'
Set-TextValue 'C50' 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
Set-TextValue 'D50' '1.004'
Set-TextValue 'E50' '  +0.41%  '
Set-TextValue 'B51' 'Decentraland'
Set-TextValue 'C51' 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue 'D51' '0.3405'
Set-TextValue 'E51' '  +1.69%  '
